$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TFEC")

# Update row 6 data (was "Solar thermal" / RES_CWH_SOLAR / "Solar thermal heating" / "Solar")
# now becomes "Biomass stoves" / RES_CWH_BIO_001 / "Biomass stoves" / "Biomass"
$ws.Range("A6").Value = "Biomass stoves"
$ws.Range("B6").Value = "RES_CWH_BIO_001"
$ws.Range("C6").Value = "Biomass stoves"
$ws.Range("D6").Value = "Biomass"

# Update the selected cell on the sheet view
$ws.Range("G8").Select()
